$d = $word.ActiveDocument

# --- 1) Drop the trailing empty paragraph and the "test" paragraph ---------
# The document starts out as:
#   Para 1: "Definisanje opstih koraka SCRUM metodologije"  (Title style)
#   Para 2: "" (empty)
#   Para 3: "test"
# After the edit, paragraphs 2 and 3 are gone entirely, so the title
# paragraph mark is immediately followed by the section properties.
while ($d.Paragraphs.Count -gt 1) {
    $lastIndex = $d.Paragraphs.Count
    $last = $d.Paragraphs($lastIndex)
    if ($lastIndex -eq 2) {
        # Only the empty paragraph remains after the title paragraph: merge
        # its paragraph mark into the title paragraph's own mark so the
        # title paragraph (with its original identity/formatting) survives
        # as the sole remaining paragraph.
        $mergeStart = $d.Paragraphs(1).Range.End - 1
        $mergeEnd = $last.Range.End
        $d.Range($mergeStart, $mergeEnd).Delete()
    }
    else {
        # Remove the whole paragraph (its text plus its own paragraph mark).
        $last.Range.Delete()
    }
}

# --- 2) Re-author the title text as individually spell-checked runs --------
# Word's background proofer now treats each Serbian word as a separate
# checked span, wrapping every word in <w:proofErr type="spellStart"/> ...
# <w:proofErr type="spellEnd"/> and breaking the former multi-word runs into
# one run per word (plus separate runs for the connecting spaces), while the
# existing "s-caron" run (different formatting) is left as-is.
$titlePara = $d.Paragraphs(1)
$titleStart = $titlePara.Range.Start
$titleTextEnd = $titlePara.Range.End - 1

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Definisanje</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> op</w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>š</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tih</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>koraka</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> SCRUM </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>metodologije</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@

$d.Range($titleStart, $titleTextEnd).InsertXML($xml)
